# Rounds_Data.xlsx update: append the two newest rounds (Mahunga on
# 2025-10-25 and Masterton Golf Course on 2025-10-26) that were logged
# after the last save, matching the "Rounds" sheet's existing layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New round-of-golf data rows (146-226) ---
$rows = @(
  @("Mahunga",1,"Russell",6),
  @("Mahunga",2,"Russell",8),
  @("Mahunga",3,"Russell",4),
  @("Mahunga",4,"Russell",4),
  @("Mahunga",5,"Russell",8),
  @("Mahunga",6,"Russell",6),
  @("Mahunga",7,"Russell",6),
  @("Mahunga",8,"Russell",7),
  @("Mahunga",9,"Russell",4),
  @("Mahunga",10,"Russell",4),
  @("Mahunga",11,"Russell",6),
  @("Mahunga",12,"Russell",6),
  @("Mahunga",13,"Russell",5),
  @("Mahunga",14,"Russell",9),
  @("Mahunga",15,"Russell",6),
  @("Mahunga",16,"Russell",5),
  @("Mahunga",17,"Russell",4),
  @("Mahunga",18,"Russell",6),
  @("Mahunga",1,"Hayden",7),
  @("Mahunga",2,"Hayden",6),
  @("Mahunga",3,"Hayden",6),
  @("Mahunga",4,"Hayden",7),
  @("Mahunga",5,"Hayden",9),
  @("Mahunga",6,"Hayden",4),
  @("Mahunga",7,"Hayden",5),
  @("Mahunga",8,"Hayden",9),
  @("Mahunga",9,"Hayden",8),
  @("Mahunga",10,"Hayden",5),
  @("Mahunga",11,"Hayden",7),
  @("Mahunga",12,"Hayden",7),
  @("Mahunga",13,"Hayden",8),
  @("Mahunga",14,"Hayden",8),
  @("Mahunga",15,"Hayden",4),
  @("Mahunga",16,"Hayden",6),
  @("Mahunga",17,"Hayden",5),
  @("Mahunga",18,"Hayden",10),
  @("Masterton Golf Course",1,"Russell",7),
  @("Masterton Golf Course",2,"Russell",9),
  @("Masterton Golf Course",3,"Russell",7),
  @("Masterton Golf Course",4,"Russell",5),
  @("Masterton Golf Course",5,"Russell",4),
  @("Masterton Golf Course",6,"Russell",7),
  @("Masterton Golf Course",7,"Russell",6),
  @("Masterton Golf Course",8,"Russell",4),
  @("Masterton Golf Course",9,"Russell",6),
  @("Masterton Golf Course",10,"Russell",8),
  @("Masterton Golf Course",11,"Russell",7),
  @("Masterton Golf Course",12,"Russell",7),
  @("Masterton Golf Course",13,"Russell",5),
  @("Masterton Golf Course",14,"Russell",6),
  @("Masterton Golf Course",15,"Russell",8),
  @("Masterton Golf Course",16,"Russell",8),
  @("Masterton Golf Course",17,"Russell",5),
  @("Masterton Golf Course",18,"Russell",6),
  @("Masterton Golf Course",1,"Hayden",9),
  @("Masterton Golf Course",2,"Hayden",9),
  @("Masterton Golf Course",3,"Hayden",3),
  @("Masterton Golf Course",4,"Hayden",8),
  @("Masterton Golf Course",5,"Hayden",6),
  @("Masterton Golf Course",6,"Hayden",6),
  @("Masterton Golf Course",7,"Hayden",7),
  @("Masterton Golf Course",8,"Hayden",6),
  @("Masterton Golf Course",9,"Hayden",8),
  @("Masterton Golf Course",10,"Hayden",7),
  @("Masterton Golf Course",11,"Hayden",4),
  @("Masterton Golf Course",12,"Hayden",10),
  @("Masterton Golf Course",13,"Hayden",7),
  @("Masterton Golf Course",14,"Hayden",10),
  @("Masterton Golf Course",15,"Hayden",11),
  @("Masterton Golf Course",16,"Hayden",7),
  @("Masterton Golf Course",17,"Hayden",6),
  @("Masterton Golf Course",18,"Hayden",5),
  @("Masterton Golf Course",1,"Olivia",14),
  @("Masterton Golf Course",2,"Olivia",9),
  @("Masterton Golf Course",3,"Olivia",6),
  @("Masterton Golf Course",4,"Olivia",6),
  @("Masterton Golf Course",5,"Olivia",11),
  @("Masterton Golf Course",6,"Olivia",9),
  @("Masterton Golf Course",7,"Olivia",10),
  @("Masterton Golf Course",8,"Olivia",8),
  @("Masterton Golf Course",9,"Olivia",9)
)

$arr = New-Object 'object[,]' $rows.Count,4
for ($i = 0; $i -lt $rows.Count; $i++) {
    for ($j = 0; $j -lt 4; $j++) {
        $arr[$i, $j] = $rows[$i][$j]
    }
}
$ws.Range("A146:D226").Value = $arr

# --- Date column (E) formatting + values ---
# Copy the existing date-style format down over the whole newly used range
# (including the trailing formatted-but-empty rows 227:235), then fill in
# the actual date serials for the data rows.
$ws.Range("E2").Copy()
$ws.Range("E146:E235").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("E146:E181").Value = 45955
$ws.Range("E182:E226").Value = 45956

# --- Column A width (widened to fit the new "Masterton Golf Course" text) ---
$ws.Columns.Item(1).ColumnWidth = 19.75

# --- Sheet view: leave the selection where the user was last working ---
$ws.Range("K219").Select() | Out-Null
